# This workbook's data rows (5-61) were re-synced from the source system and
# ended up re-sorted: every data row now carries the payload that used to
# belong to a *different* row, while the row-number (and therefore any
# formatting keyed off it) stays put. The mapping below says, for each target
# row, which row currently (pre-edit) holds the data that must end up there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetToSource = @{
    5 = 41;
    6 = 5;
    7 = 42;
    8 = 43;
    9 = 44;
    10 = 45;
    11 = 6;
    12 = 7;
    13 = 8;
    14 = 9;
    15 = 46;
    16 = 10;
    17 = 11;
    18 = 47;
    19 = 48;
    20 = 49;
    21 = 50;
    22 = 12;
    23 = 13;
    24 = 51;
    25 = 14;
    26 = 15;
    27 = 52;
    28 = 16;
    29 = 17;
    30 = 18;
    31 = 19;
    32 = 53;
    33 = 54;
    34 = 20;
    35 = 21;
    36 = 22;
    37 = 55;
    38 = 56;
    39 = 23;
    40 = 24;
    41 = 25;
    42 = 57;
    43 = 26;
    44 = 27;
    45 = 58;
    46 = 59;
    47 = 28;
    48 = 29;
    49 = 30;
    50 = 31;
    51 = 32;
    52 = 60;
    53 = 33;
    54 = 34;
    55 = 35;
    56 = 61;
    57 = 36;
    58 = 37;
    59 = 38;
    60 = 39;
    61 = 40
}

# Snapshot every source row (columns A:AY) BEFORE any writes happen, since the
# mapping is a permutation and rows are both sources and destinations.
$snapshot = @{}
foreach ($row in $targetToSource.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rng = $ws.Range("A$row`:AY$row")
        $snapshot[$row] = $rng.Value2
    }
}

# Now write each target row from the snapshot of its mapped source row.
foreach ($row in ($targetToSource.Keys | Sort-Object)) {
    $src = $targetToSource[$row]
    $destRange = $ws.Range("A$row`:AY$row")
    $destRange.Value2 = $snapshot[$src]
}
